# CVX.xlsx model update — "Add files via upload"
#
# Core change: the long-term FCF-margin growth assumption (row 28, Sheet2)
# is revised down from 3.5%/yr to 2.0%/yr, starting at E28 and filling
# across the projection (F28:N28 share that formula). This ripples through
# FCF (row 32), FCF/NI (row 33) and the NPV/valuation block (Q23:Q25).
#
# Along with that, a couple of label/total cells (P23 "NPV", and the
# "FCF Margin" row-25 helper cells A25:D25) have bold formatting removed,
# and row 14 (E14:N14) is re-entered as one fill-across formula (same
# values, now a single shared formula group). The last active cell on the
# sheet is left on N28 (the updated growth-rate cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 14: re-enter the fill-across formula (D34*$Q$20 … M34*$Q$20) ---
# Same values as before; Excel now stores E14:N14 as a single shared
# formula instead of ten separate ones.
$ws.Range("E14:N14").Formula = '=D34*$Q$20'

# --- Row 28: lower the post-projection FCF-margin growth rate 3.5% -> 2% ---
$ws.Range("E28").Formula = '=D28*1.02'
$ws.Range("F28:N28").Formula = '=E28*1.02'

# --- Formatting cleanup: drop bold on the NPV label and the FCF Margin
#     helper row ---
$ws.Range("P23").Font.Bold = $false
$ws.Range("A25:D25").Font.Bold = $false

# --- Leave the selection on the cell that was last edited ---
$ws.Range("N28").Select()
